{"js": "// Update the worksheet date and all two-digit division problems to the\n// new values per the target revision.\nconst replacements = [\n  [\"2023-12-23 Saturday\", \"2023-12-24 Sunday\"],\n  [\"98\u00f79=\", \"27\u00f79=\"],\n  [\"72\u00f77=\", \"14\u00f76=\"],\n  [\"59\u00f72=\", \"25\u00f73=\"],\n  [\"33\u00f74=\", \"49\u00f77=\"],\n  [\"82\u00f78=\", \"36\u00f73=\"],\n  [\"41\u00f75=\", \"84\u00f76=\"],\n  [\"90\u00f76=\", \"44\u00f77=\"],\n  [\"99\u00f75=\", \"89\u00f73=\"],\n  [\"74\u00f77=\", \"55\u00f76=\"],\n  [\"17\u00f72=\", \"24\u00f73=\"],\n  [\"78\u00f75=\", \"83\u00f74=\"],\n  [\"70\u00f78=\", \"83\u00f72=\"],\n  [\"44\u00f72=\", \"17\u00f78=\"],\n  [\"80\u00f76=\", \"77\u00f73=\"],\n  [\"18\u00f72=\", \"79\u00f78=\"],\n  [\"57\u00f74=\", \"29\u00f75=\"],\n  [\"54\u00f77=\", \"81\u00f75=\"],\n  [\"13\u00f79=\", \"50\u00f79=\"],\n  [\"96\u00f73=\", \"25\u00f78=\"],\n  [\"25\u00f72=\", \"26\u00f78=\"],\n  [\"59\u00f78=\", \"36\u00f78=\"],\n  [\"65\u00f77=\", \"57\u00f78=\"],\n  [\"68\u00f72=\", \"30\u00f75=\"],\n  [\"43\u00f76=\", \"55\u00f76=\"],\n  [\"95\u00f77=\", \"21\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all two-digit division problems to the\n# new values per the target revision.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-12-23 Saturday\", \"2023-12-24 Sunday\"),\n    @(\"98\u00f79=\", \"27\u00f79=\"),\n    @(\"72\u00f77=\", \"14\u00f76=\"),\n    @(\"59\u00f72=\", \"25\u00f73=\"),\n    @(\"33\u00f74=\", \"49\u00f77=\"),\n    @(\"82\u00f78=\", \"36\u00f73=\"),\n    @(\"41\u00f75=\", \"84\u00f76=\"),\n    @(\"90\u00f76=\", \"44\u00f77=\"),\n    @(\"99\u00f75=\", \"89\u00f73=\"),\n    @(\"74\u00f77=\", \"55\u00f76=\"),\n    @(\"17\u00f72=\", \"24\u00f73=\"),\n    @(\"78\u00f75=\", \"83\u00f74=\"),\n    @(\"70\u00f78=\", \"83\u00f72=\"),\n    @(\"44\u00f72=\", \"17\u00f78=\"),\n    @(\"80\u00f76=\", \"77\u00f73=\"),\n    @(\"18\u00f72=\", \"79\u00f78=\"),\n    @(\"57\u00f74=\", \"29\u00f75=\"),\n    @(\"54\u00f77=\", \"81\u00f75=\"),\n    @(\"13\u00f79=\", \"50\u00f79=\"),\n    @(\"96\u00f73=\", \"25\u00f78=\"),\n    @(\"25\u00f72=\", \"26\u00f78=\"),\n    @(\"59\u00f78=\", \"36\u00f78=\"),\n    @(\"65\u00f77=\", \"57\u00f78=\"),\n    @(\"68\u00f72=\", \"30\u00f75=\"),\n    @(\"43\u00f76=\", \"55\u00f76=\"),\n    @(\"95\u00f77=\", \"21\u00f75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
